$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.176
$ws.Range("C2").Value = 0.592
$ws.Range("J2").Value = 0.008
$ws.Range("P2").Value = 0.132
$ws.Range("S2").Value = 0.092
$ws.Range("B3").Value = 0.006451612903225806
$ws.Range("C3").Value = 0.03870967741935484
$ws.Range("P3").Value = 0.8193548387096774
$ws.Range("S3").Value = 0.1354838709677419
$ws.Range("J4").Value = 0.1428571428571428
$ws.Range("P4").Value = 0.5142857142857142
$ws.Range("S4").Value = 0.3428571428571429
$ws.Range("J5").Value = 0.5
$ws.Range("P5").Value = 0.5
$ws.Range("B6").Value = 0.08900523560209424
$ws.Range("D6").Value = 0.02617801047120419
$ws.Range("E6").Value = 0.005235602094240838
$ws.Range("F6").Value = 0.06282722513089005
$ws.Range("J6").Value = 0.2094240837696335
$ws.Range("O6").Value = 0.005235602094240838
$ws.Range("Q6").Value = 0.1047120418848168
$ws.Range("R6").Value = 0.1047120418848168
$ws.Range("S6").Value = 0.3926701570680629
$ws.Range("B7").Value = 0.1290322580645161
$ws.Range("D7").Value = 0.01612903225806452
$ws.Range("F7").Value = 0.06451612903225806
$ws.Range("J7").Value = 0.1182795698924731
$ws.Range("O7").Value = 0.005376344086021506
$ws.Range("Q7").Value = 0.2096774193548387
$ws.Range("R7").Value = 0.06989247311827956
$ws.Range("S7").Value = 0.3870967741935484
$ws.Range("B8").Value = 0.07363420427553444
$ws.Range("D8").Value = 0.01900237529691211
$ws.Range("F8").Value = 0.05938242280285035
$ws.Range("J8").Value = 0.1187648456057007
$ws.Range("O8").Value = 0.007125890736342043
$ws.Range("Q8").Value = 0.1496437054631829
$ws.Range("R8").Value = 0.1235154394299287
$ws.Range("S8").Value = 0.4489311163895487
$ws.Range("B9").Value = 0.1169590643274854
$ws.Range("D9").Value = 0.01754385964912281
$ws.Range("E9").Value = 0.005847953216374269
$ws.Range("F9").Value = 0.05263157894736842
$ws.Range("J9").Value = 0.1052631578947368
$ws.Range("Q9").Value = 0.2046783625730994
$ws.Range("R9").Value = 0.05263157894736842
$ws.Range("S9").Value = 0.4444444444444444
$ws.Range("B10").Value = 0.09632224168126094
$ws.Range("D10").Value = 0.0148861646234676
$ws.Range("F10").Value = 0.06830122591943957
$ws.Range("J10").Value = 0.1348511383537653
$ws.Range("O10").Value = 0.003502626970227671
$ws.Range("Q10").Value = 0.2066549912434326
$ws.Range("R10").Value = 0.07793345008756568
$ws.Range("S10").Value = 0.3975481611208406
$ws.Range("G11").Value = 0.1596091205211726
$ws.Range("J11").Value = 0.07817589576547231
$ws.Range("K11").Value = 0.2214983713355049
$ws.Range("L11").Value = 0.5211726384364821
$ws.Range("S11").Value = 0.01954397394136808
$ws.Range("G12").Value = 0.7391304347826086
$ws.Range("J12").Value = 0.2049689440993789
$ws.Range("K12").Value = 0.006211180124223602
$ws.Range("L12").Value = 0.006211180124223602
$ws.Range("S12").Value = 0.04347826086956522
$ws.Range("G13").Value = 0.6923076923076923
$ws.Range("J13").Value = 0.2051282051282051
$ws.Range("S13").Value = 0.1025641025641026
$ws.Range("F15").Value = 0.01796407185628742
$ws.Range("H15").Value = 0.1976047904191617
$ws.Range("I15").Value = 0.07784431137724551
$ws.Range("J15").Value = 0.407185628742515
$ws.Range("K15").Value = 0.0658682634730539
$ws.Range("M15").Value = 0.01796407185628742
$ws.Range("O15").Value = 0.07784431137724551
$ws.Range("S15").Value = 0.1377245508982036
$ws.Range("F16").Value = 0.02857142857142857
$ws.Range("H16").Value = 0.1485714285714286
$ws.Range("I16").Value = 0.08
$ws.Range("J16").Value = 0.3942857142857143
$ws.Range("K16").Value = 0.08571428571428572
$ws.Range("M16").Value = 0.005714285714285714
$ws.Range("O16").Value = 0.05714285714285714
$ws.Range("S16").Value = 0.2
$ws.Range("F17").Value = 0.0102827763496144
$ws.Range("H17").Value = 0.1979434447300771
$ws.Range("I17").Value = 0.08997429305912596
$ws.Range("J17").Value = 0.3984575835475578
$ws.Range("K17").Value = 0.1053984575835476
$ws.Range("M17").Value = 0.01542416452442159
$ws.Range("O17").Value = 0.05655526992287917
$ws.Range("S17").Value = 0.1259640102827763
$ws.Range("F18").Value = 0.01630434782608696
$ws.Range("H18").Value = 0.2173913043478261
$ws.Range("I18").Value = 0.05978260869565218
$ws.Range("J18").Value = 0.3858695652173913
$ws.Range("K18").Value = 0.1141304347826087
$ws.Range("M18").Value = 0.01630434782608696
$ws.Range("O18").Value = 0.05434782608695652
$ws.Range("S18").Value = 0.1358695652173913
$ws.Range("F19").Value = 0.01137286758732738
$ws.Range("H19").Value = 0.1990251827782291
$ws.Range("I19").Value = 0.08042242079610074
$ws.Range("J19").Value = 0.3549959382615759
$ws.Range("K19").Value = 0.1186027619821284
$ws.Range("M19").Value = 0.02193338748984565
$ws.Range("O19").Value = 0.06904955320877336
$ws.Range("S19").Value = 0.1445978878960195
